# Updated symbol list on Tue Jan 17 06:32:09 UTC 2023 with GitHub Actions
# Applies updated Price (D) and Volume(1h) (E) text values to the cryptos worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "299.04"
Set-TextValue $ws.Range("E2") "-1.67%"
Set-TextValue $ws.Range("E3") "-1.33%"
Set-TextValue $ws.Range("D4") "5.114"
Set-TextValue $ws.Range("E4") "-2.11%"
Set-TextValue $ws.Range("D5") "0.07915"
Set-TextValue $ws.Range("E5") "1.42%"
Set-TextValue $ws.Range("D6") "2.350"
Set-TextValue $ws.Range("E6") "2.91%"
Set-TextValue $ws.Range("D7") "7.804"
Set-TextValue $ws.Range("E7") "-2.47%"
Set-TextValue $ws.Range("D8") "3.865"
Set-TextValue $ws.Range("E8") "-0.20%"
Set-TextValue $ws.Range("D9") "0.9234"
Set-TextValue $ws.Range("E9") "1.10%"
Set-TextValue $ws.Range("D10") "0.1748"
Set-TextValue $ws.Range("E10") "0.34%"
Set-TextValue $ws.Range("D11") "0.07561"
Set-TextValue $ws.Range("E11") "2.48%"
Set-TextValue $ws.Range("D12") "0.09329"
Set-TextValue $ws.Range("E12") "13.56%"
Set-TextValue $ws.Range("D13") "0.03010"
Set-TextValue $ws.Range("E13") "-1.07%"
Set-TextValue $ws.Range("E14") "0.59%"
Set-TextValue $ws.Range("D15") "0.001509"
Set-TextValue $ws.Range("E15") "-1.04%"
Set-TextValue $ws.Range("D16") "0.005949"
Set-TextValue $ws.Range("E16") "-3.73%"
Set-TextValue $ws.Range("D17") "3.475"
Set-TextValue $ws.Range("E17") "-0.71%"
Set-TextValue $ws.Range("D18") "2.265"
Set-TextValue $ws.Range("E18") "1.38%"
Set-TextValue $ws.Range("E19") "0.23%"
Set-TextValue $ws.Range("E20") "-1.54%"
Set-TextValue $ws.Range("D21") "3.902"
Set-TextValue $ws.Range("D22") "0.1701"
Set-TextValue $ws.Range("E22") "8.73%"
Set-TextValue $ws.Range("E23") "-0.35%"
Set-TextValue $ws.Range("D24") "0.001252"
Set-TextValue $ws.Range("E24") "-0.81%"
Set-TextValue $ws.Range("D25") "0.004470"
Set-TextValue $ws.Range("E25") "-1.41%"
Set-TextValue $ws.Range("D26") "0.0001200"
Set-TextValue $ws.Range("E26") "-11.02%"
Set-TextValue $ws.Range("D27") "0.0003398"
Set-TextValue $ws.Range("E27") "24.03%"
Set-TextValue $ws.Range("D39") "0.01744"
Set-TextValue $ws.Range("E39") "-2.48%"
Set-TextValue $ws.Range("D40") "0.04612"
Set-TextValue $ws.Range("E40") "0.43%"
Set-TextValue $ws.Range("D41") "0.006967"
Set-TextValue $ws.Range("E41") "-4.29%"
Set-TextValue $ws.Range("E42") "-0.23%"
Set-TextValue $ws.Range("D43") "0.002220"
Set-TextValue $ws.Range("E43") "-0.80%"
Set-TextValue $ws.Range("E44") "-5.55%"
Set-TextValue $ws.Range("D45") "0.00006284"
Set-TextValue $ws.Range("E45") "-2.54%"
Set-TextValue $ws.Range("E46") "0.11%"
Set-TextValue $ws.Range("D47") "0.007982"
Set-TextValue $ws.Range("E47") "-19.30%"
Set-TextValue $ws.Range("D48") "1.156"
Set-TextValue $ws.Range("E48") "40.82%"
Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("E49") "0.11%"
Set-TextValue $ws.Range("D50") "0.0002000"
Set-TextValue $ws.Range("E50") "0.11%"
